# Auto-generated edit script replicating the XML diff for Halicarnassus_Profits.xlsx
# Applies literal value changes / cell clears / one cell addition across 5 sheets
# (ALC, ARM, CUL, GSM, LTW, WVR) matching the scheduled market-data refresh commit.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets("ALC")
# Row 2
$ws.Range("H2").Value = 173.11765
$ws.Range("I2").Value = 179.41667
$ws.Range("K2").Value = 179.41667
$ws.Range("M2").Value = -66.41667000000001
# Row 43
$ws.Range("H43").Value = 4772779.5
$ws.Range("J43").Value = 23272.25
$ws.Range("L43").Value = 23272.25
$ws.Range("N43").Value = -23410.25
# Row 74
$ws.Range("H74").Value = 6662
$ws.Range("I74").Value = 6662
$ws.Range("K74").Value = 6662
$ws.Range("M74").Value = -5726
# Row 77
$ws.Range("H77").Value = 6662
$ws.Range("I77").Value = 6662
$ws.Range("K77").Value = 33310
$ws.Range("M77").Value = -28630
# Row 80
$ws.Range("H80").Value = 742.2222
$ws.Range("I80").Value = 808
$ws.Range("J80").Value = 660
$ws.Range("K80").Value = 2424
$ws.Range("L80").Value = 1980
$ws.Range("M80").Value = -1426
$ws.Range("N80").Value = -3976
# Row 83
$ws.Range("H83").Value = 742.2222
$ws.Range("I83").Value = 808
$ws.Range("J83").Value = 660
$ws.Range("K83").Value = 7272
$ws.Range("L83").Value = 5940
$ws.Range("M83").Value = -2280
$ws.Range("N83").Value = -15924
# Row 137
$ws.Range("H137").Value = 2873.4482
$ws.Range("I137").Value = 2131.2856
$ws.Range("K137").Value = 6393.8568
$ws.Range("M137").Value = -3843.8568

# ---- ARM sheet ----
$ws = $wb.Worksheets("ARM")
# Row 74
$ws.Range("H74").Value = 2273.6538
$ws.Range("I74").Value = 1996.1364
$ws.Range("J74").Value = 3800
$ws.Range("K74").Value = 1996.1364
$ws.Range("L74").Value = 3800
$ws.Range("M74").Value = -1122.1364
$ws.Range("N74").Value = -5548
# Row 77
$ws.Range("H77").Value = 2273.6538
$ws.Range("I77").Value = 1996.1364
$ws.Range("J77").Value = 3800
$ws.Range("K77").Value = 9980.682000000001
$ws.Range("L77").Value = 19000
$ws.Range("M77").Value = -5612.682000000001
$ws.Range("N77").Value = -27736

# ---- CUL sheet ----
$ws = $wb.Worksheets("CUL")
# Row 107
$ws.Range("H107").Value = 461.66666
$ws.Range("J107").Value = 643.5
$ws.Range("L107").Value = 1930.5
$ws.Range("N107").Value = -5770.5
# Row 120
$ws.Range("H120").ClearContents()
$ws.Range("I120").ClearContents()
$ws.Range("J120").ClearContents()
$ws.Range("K120").ClearContents()
$ws.Range("L120").ClearContents()
$ws.Range("M120").ClearContents()
# Row 121
$ws.Range("H121").ClearContents()
$ws.Range("I121").ClearContents()
$ws.Range("J121").ClearContents()
$ws.Range("K121").ClearContents()
$ws.Range("L121").ClearContents()
$ws.Range("M121").ClearContents()
$ws.Range("N121").ClearContents()
# Row 122
$ws.Range("H122").ClearContents()
$ws.Range("I122").ClearContents()
$ws.Range("J122").ClearContents()
$ws.Range("K122").ClearContents()
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
# Row 123
$ws.Range("H123").ClearContents()
$ws.Range("I123").ClearContents()
$ws.Range("J123").ClearContents()
$ws.Range("K123").ClearContents()
$ws.Range("L123").ClearContents()
$ws.Range("M123").ClearContents()
# Row 124
$ws.Range("H124").ClearContents()
$ws.Range("I124").ClearContents()
$ws.Range("J124").ClearContents()
$ws.Range("K124").ClearContents()
$ws.Range("L124").ClearContents()
$ws.Range("M124").ClearContents()
# Row 125
$ws.Range("H125").ClearContents()
$ws.Range("I125").ClearContents()
$ws.Range("J125").ClearContents()
$ws.Range("K125").ClearContents()
$ws.Range("L125").ClearContents()
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
# Row 126
$ws.Range("H126").ClearContents()
$ws.Range("I126").ClearContents()
$ws.Range("J126").ClearContents()
$ws.Range("K126").ClearContents()
$ws.Range("L126").ClearContents()
# Row 127
$ws.Range("H127").ClearContents()
$ws.Range("I127").ClearContents()
$ws.Range("J127").ClearContents()
$ws.Range("K127").ClearContents()
$ws.Range("L127").ClearContents()
# Row 128
$ws.Range("H128").ClearContents()
$ws.Range("I128").ClearContents()
$ws.Range("J128").ClearContents()
$ws.Range("K128").ClearContents()
$ws.Range("L128").ClearContents()
$ws.Range("M128").ClearContents()
# Row 129
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()
# Row 130
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("M130").ClearContents()
$ws.Range("N130").ClearContents()
# Row 131
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("M131").ClearContents()
$ws.Range("N131").ClearContents()
# Row 132
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
# Row 133
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
# Row 134
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
# Row 136
$ws.Range("H136").ClearContents()
$ws.Range("I136").ClearContents()
$ws.Range("J136").ClearContents()
$ws.Range("K136").ClearContents()
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()
# Row 137
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("M137").ClearContents()
# Row 138
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()
# Row 139
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()
# Row 140
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("M140").ClearContents()
$ws.Range("N140").ClearContents()
# Row 141
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

# ---- GSM sheet ----
$ws = $wb.Worksheets("GSM")
# Row 18
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
# Row 46
$ws.Range("H46").Value = 37935.332
$ws.Range("I46").Value = 7520.5
$ws.Range("K46").Value = 7520.5
$ws.Range("M46").Value = -7364.5
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("N67").ClearContents()

# ---- LTW sheet ----
$ws = $wb.Worksheets("LTW")
# Row 122
$ws.Range("H122").Value = 2897.5
$ws.Range("J122").Value = 2895
$ws.Range("L122").Value = 8685
$ws.Range("N122").Value = -13585

# ---- WVR sheet ----
$ws = $wb.Worksheets("WVR")
# Row 41
$ws.Range("H41").Value = 23101.727
$ws.Range("I41").Value = 20705.6
$ws.Range("J41").Value = 25098.5
$ws.Range("K41").Value = 20705.6
$ws.Range("L41").Value = 25098.5
$ws.Range("M41").Value = -20315.6
$ws.Range("N41").Value = -25878.5
# Row 100
$ws.Range("H100").Value = 892.2308
$ws.Range("I100").Value = 799.9167
$ws.Range("K100").Value = 1599.8334
$ws.Range("M100").Value = -1058.8334
# Row 136
$ws.Range("H136").Value = 3292.6333
$ws.Range("I136").Value = 2269.9473
$ws.Range("J136").Value = 5059.091
$ws.Range("K136").Value = 6809.841899999999
$ws.Range("L136").Value = 15177.273
$ws.Range("M136").Value = -4259.841899999999
$ws.Range("N136").Value = -20277.273
